# Fixed a bug in GameScenePoolEx
# The data rows (2-21) on the active sheet were reordered / corrected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(701,  3, 90, 45, 97, 15),
    @(801,  3, 67, 65, 52, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(1203, 3, 15, 15, 15, 15),
    @(901, 16, 15, 45, 60, 60),
    @(301,  6, 45, 30, 60, 45),
    @(401,  9, 48, 67, 75, 45),
    @(601,  9, 60, 67, 60, 42),
    @(201,  9, 30, 15, 45, 30),
    @(101,  9, 30, 15, 60, 15),
    @(902,  1,  0,  0,  0,  0),
    @(1001,18, 30, 75, 60, 72),
    @(501,  9, 52, 30, 75, 45),
    @(802,  0,  4,  5,  4,  0),
    @(502,  0,  4,  0,  0,  0),
    @(1,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3),
    @(1101, 0, 15, 30, 30,  0),
    @(2,    0,  2,  2,  2,  2)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
